$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '57.793.81'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +2.18%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.056.10'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +2.23%  '

$ws.Range("E4").Value = '  +0.02%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '524.65'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +5.50%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '142.40'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +4.93%  '

$ws.Range("E7").Value = '  +0.03%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.447'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +4.72%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '7.63'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +4.99%  '

$ws.Range("E10").Value = '  +7.64%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.372'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +4.93%  '

$ws.Range("E12").Value = '  +2.13%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '3.578.48'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +2.27%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '26.88'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +7.27%  '

$ws.Range("E15").Value = '  +17.04%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '57.755.68'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +2.21%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '6.25'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +6.64%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.045.05'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.97%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '13.07'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +5.19%  '

$ws.Range("E20").Value = '  +5.09%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '339.08'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +4.03%  '

$ws.Range("E22").Value = '  +0.13%  '

$ws.Range("E23").Value = '  +7.07%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '64.96'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +5.82%  '

$ws.Range("E25").Value = '  +5.83%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.0₃0979'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +8.45%  '

$ws.Range("E27").Value = '  +0.37%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '6.93'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +6.24%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '7.38'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +10.34%  '

$ws.Range("E30").Value = '  +7.15%  '

$ws.Range("E31").Value = '  +4.71%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '21.15'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +4.27%  '

$ws.Range("E33").Value = '  +6.18%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '156.70'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +0.62%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '5.98'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +6.23%  '

$ws.Range("E36").Value = '  +2.40%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '26.19'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +11.74%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.0705'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +2.88%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '3.091.46'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.34%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '37.72'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +3.05%  '

$ws.Range("E41").Value = '  +8.75%  '

$ws.Range("E42").Value = '  +0.02%  '

$ws.Range("E43").Value = '  +5.28%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.664'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +3.90%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.333.97'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +4.73%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.04'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +3.99%  '

$ws.Range("E47").Value = '  +4.12%  '

$ws.Range("E48").Value = '  +4.82%  '

$ws.Range("E49").Value = '  +4.04%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '20.26'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +5.87%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0898'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +5.70%  '
